# Consolidate the evaluation sheets under clearer, descriptive names
# ("a colocar tudo no mesmo excel").
#
# Note: Excel enforces a hard 31-character limit on worksheet names, so the
# desired names are truncated to the maximum Excel allows (real Excel would
# refuse anything longer, raising "Invalid name for a sheet").
$wb = $excel.ActiveWorkbook

$vpe = $wb.Worksheets.Item("VPE")
$vpe.Name = "Avaliacao-Vice-Presidente-Exter"

$rh = $wb.Worksheets.Item("RH")
$rh.Name = "Avaliacao-Membro-RH"
